# Updated cryptos list with GitHub Actions
# Applies the latest price / 1h-volume-change snapshot to the crypto
# tracker sheet, and fixes the Bittensor / RenderToken row ordering.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.084.48"
$ws.Range("E2").Value = "  -0.09%  "
$ws.Range("D3").Value = "2.758.05"
$ws.Range("E3").Value = "  +0.55%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "'575.84"
$ws.Range("E5").Value = "  -0.77%  "
$ws.Range("D6").Value = "'158.95"
$ws.Range("E6").Value = "  +0.54%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").Value = "'0.604"
$ws.Range("E8").Value = "  -3.26%  "
$ws.Range("E9").Value = "  -1.84%  "
$ws.Range("D10").Value = "'0.167"
$ws.Range("E10").Value = "  +4.60%  "
$ws.Range("D12").Value = "'5.70"
$ws.Range("E12").Value = "  -15.89%  "
$ws.Range("D13").Value = "3.249.01"
$ws.Range("E13").Value = "  +0.25%  "
$ws.Range("D14").Value = "'26.91"
$ws.Range("E14").Value = "  -1.64%  "
$ws.Range("D15").Value = "63.723.33"
$ws.Range("E15").Value = "  -0.45%  "
$ws.Range("E16").Value = "  -2.65%  "
$ws.Range("D17").Value = "2.762.92"
$ws.Range("E17").Value = "  +0.38%  "
$ws.Range("D18").Value = "'12.13"
$ws.Range("E18").Value = "  -0.03%  "
$ws.Range("E19").Value = "  -1.69%  "
$ws.Range("D20").Value = "'357.68"
$ws.Range("E20").Value = "  -1.85%  "
$ws.Range("D21").Value = "'6.75"
$ws.Range("E21").Value = "  -4.16%  "
$ws.Range("D22").Value = "'0.999"
$ws.Range("E22").Value = "  +0.16%  "
$ws.Range("D23").Value = "'0.536"
$ws.Range("E23").Value = "  -1.35%  "
$ws.Range("D24").Value = "'65.59"
$ws.Range("E24").Value = "  -1.89%  "
$ws.Range("E25").Value = "  -1.23%  "
$ws.Range("D26").Value = "'8.60"
$ws.Range("E26").Value = "  -0.75%  "
$ws.Range("E27").Value = "  +0.23%  "
$ws.Range("D28").Value = "0.0₃0909"
$ws.Range("E28").Value = "  -1.37%  "
$ws.Range("D29").Value = "'7.28"
$ws.Range("E29").Value = "  +0.05%  "
$ws.Range("E30").Value = "  -3.16%  "
$ws.Range("E31").Value = "  -1.02%  "
$ws.Range("D32").Value = "'169.35"
$ws.Range("E32").Value = "  -2.73%  "
$ws.Range("D33").Value = "'20.27"
$ws.Range("E33").Value = "  -1.92%  "
$ws.Range("D34").Value = "'4.93"
$ws.Range("E34").Value = "  -0.61%  "
$ws.Range("D35").Value = "'1.49"
$ws.Range("E35").Value = "  +1.97%  "
$ws.Range("E36").Value = "  -0.04%  "
$ws.Range("E37").Value = "  -0.53%  "
$ws.Range("D38").Value = "'1.00"
$ws.Range("E38").Value = "  -0.60%  "
$ws.Range("B39").Value = "Bittensor"
$ws.Range("C39").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D39").Value = "'346.05"
$ws.Range("E39").Value = "  +2.01%  "
$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D40").Value = "'6.33"
$ws.Range("E40").Value = "  +2.41%  "
$ws.Range("E41").Value = "  -2.08%  "
$ws.Range("D42").Value = "'39.22"
$ws.Range("E42").Value = "  -0.83%  "
$ws.Range("E43").Value = "  -2.04%  "
$ws.Range("D44").Value = "'21.86"
$ws.Range("E44").Value = "  -3.24%  "
$ws.Range("D45").Value = "'0.0591"
$ws.Range("E45").Value = "  -2.32%  "
$ws.Range("D46").Value = "'0.0255"
$ws.Range("E46").Value = "  -1.72%  "
$ws.Range("D47").Value = "'0.632"
$ws.Range("E47").Value = "  -2.21%  "
$ws.Range("E48").Value = "  -0.49%  "
$ws.Range("D49").Value = "'135.55"
$ws.Range("E49").Value = "  -1.75%  "
$ws.Range("D50").Value = "'0.998"
$ws.Range("E50").Value = "  -0.06%  "
$ws.Range("D51").Value = "'11.04"
$ws.Range("E51").Value = "  -0.11%  "

